$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 4848
$ws.Range("K3").Value = 4980
$ws.Range("K4").Value = 1033
$ws.Range("K6").Value = 5598
$ws.Range("K7").Value = 16815

# By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K7").Value = 495
$ws.Range("K8").Value = 1125
$ws.Range("K10").Value = 92
$ws.Range("K11").Value = 329
$ws.Range("K13").Value = 16
$ws.Range("K18").Value = 113
$ws.Range("K19").Value = 510
$ws.Range("K20").Value = 388
$ws.Range("K23").Value = 169
$ws.Range("K27").Value = 154
$ws.Range("K29").Value = 897
$ws.Range("K30").Value = 64
$ws.Range("K33").Value = 710
$ws.Range("K36").Value = 217
$ws.Range("K37").Value = 568
$ws.Range("K42").Value = 624
$ws.Range("K43").Value = 148
$ws.Range("K44").Value = 147
$ws.Range("K49").Value = 94
$ws.Range("K51").Value = 213
$ws.Range("K54").Value = 329
$ws.Range("K55").Value = 192
$ws.Range("K57").Value = 61
$ws.Range("K63").Value = 55
$ws.Range("K64").Value = 105
$ws.Range("K65").Value = 380
$ws.Range("K67").Value = 644
$ws.Range("K68").Value = 45
$ws.Range("K72").Value = 78
$ws.Range("K76").Value = 230
$ws.Range("K77").Value = 122
$ws.Range("K79").Value = 412
$ws.Range("K80").Value = 59
$ws.Range("K83").Value = 369
$ws.Range("K84").Value = 125
$ws.Range("K85").Value = 769
$ws.Range("K86").Value = 113
$ws.Range("K89").Value = 240
$ws.Range("K91").Value = 182
$ws.Range("K94").Value = 220
$ws.Range("K95").Value = 294
$ws.Range("K96").Value = 182
$ws.Range("K98").Value = 82
$ws.Range("K99").Value = 288
$ws.Range("K101").Value = 16815

# West Ridge
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K4").Value = 9
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 182

# Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 173
$ws.Range("K3").Value = 157
$ws.Range("K7").Value = 495

# Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 110
$ws.Range("K3").Value = 86
$ws.Range("K7").Value = 329

# Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 64
$ws.Range("K7").Value = 240

# South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K3").Value = 260
$ws.Range("K4").Value = 46
$ws.Range("K6").Value = 182
$ws.Range("K7").Value = 769

# Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 315
$ws.Range("K3").Value = 336
$ws.Range("K6").Value = 379
$ws.Range("K7").Value = 1125

# South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K3").Value = 134
$ws.Range("K6").Value = 85
$ws.Range("K7").Value = 369

# Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K6").Value = 203
$ws.Range("K7").Value = 710

# West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 98
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 294

# Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 161
$ws.Range("K6").Value = 169
$ws.Range("K7").Value = 568

# New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 95
$ws.Range("K7").Value = 380

# Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 74
$ws.Range("K3").Value = 118
$ws.Range("K7").Value = 288

# Fuller Park
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 64

# North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 225
$ws.Range("K7").Value = 644

# South Deering
$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K3").Value = 49
$ws.Range("K7").Value = 125

# Lincoln Park
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K2").Value = 17
$ws.Range("K7").Value = 94

# Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K4").Value = 18
$ws.Range("K6").Value = 173
$ws.Range("K7").Value = 329

# Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 260
$ws.Range("K3").Value = 321
$ws.Range("K6").Value = 247
$ws.Range("K7").Value = 897

# Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 152
$ws.Range("K3").Value = 160
$ws.Range("K4").Value = 22
$ws.Range("K7").Value = 510

# Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 147

# River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 126
$ws.Range("K7").Value = 230

# Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 166
$ws.Range("K4").Value = 23
$ws.Range("K6").Value = 237
$ws.Range("K7").Value = 624

# Boystown
$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("K5").Value = 7
$ws.Range("K6").Value = 16

# Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 92

# Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 192

# Douglas
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 169

# Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 47
$ws.Range("K7").Value = 182

# Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 136
$ws.Range("K3").Value = 131
$ws.Range("K7").Value = 412

# Near South Side
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 105

# Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 130
$ws.Range("K3").Value = 124
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 388

# Calumet Heights
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 113

# Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 85
$ws.Range("K3").Value = 60
$ws.Range("K7").Value = 217

# West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 220

# Wicker Park
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 82

# Edgewater
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 154

# Streeterville
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 113

# Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 213

# North Park
$ws = $wb.Worksheets.Item('North Park')
$ws.Range("K2").Value = 19
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 45

# Mckinley Park
$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 61

# Hyde Park
$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K3").Value = 40
$ws.Range("K7").Value = 148

# Old Town
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 78

# Riverdale
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K3").Value = 49
$ws.Range("K7").Value = 122

# Rush & Division
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 59
